$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.557.93"
$ws.Range("E2").Value = "  +0.92%  "

$ws.Range("D3").Value = "3.728.19"
$ws.Range("E3").Value = "  -1.02%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.67"
$ws.Range("E5").Value = "  +0.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.71"
$ws.Range("E6").Value = "  -4.75%  "

$ws.Range("D7").Value = "3.727.41"
$ws.Range("E7").Value = "  -1.07%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  +1.32%  "

$ws.Range("E10").Value = "  +2.55%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("E11").Value = "  +2.96%  "

$ws.Range("E12").Value = "  -1.21%  "

$ws.Range("E13").Value = "  -1.22%  "

$ws.Range("E14").Value = "  -0.81%  "

$ws.Range("D15").Value = "4.354.36"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").Value = "3.725.31"
$ws.Range("E16").Value = "  -0.83%  "

$ws.Range("D17").Value = "68.583.92"
$ws.Range("E17").Value = "  +1.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.25"
$ws.Range("E18").Value = "  +0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.116"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.24"
$ws.Range("E20").Value = "  +3.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "495.13"
$ws.Range("E21").Value = "  +0.86%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.17"
$ws.Range("E22").Value = "  +12.34%  "

$ws.Range("E23").Value = "  -2.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.91"
$ws.Range("E24").Value = "  -0.55%  "

$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000142"
$ws.Range("E25").Value = "  -5.05%  "

$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  -2.18%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.41"
$ws.Range("E27").Value = "  +1.01%  "

$ws.Range("E28").Value = "  -1.14%  "

$ws.Range("E29").Value = "  +0.07%  "

$ws.Range("E30").Value = "  -0.72%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.46"
$ws.Range("E31").Value = "  +0.68%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.91"
$ws.Range("E32").Value = "  +2.82%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.59"
$ws.Range("E33").Value = "  -2.14%  "

$ws.Range("D34").Value = "3.875.30"
$ws.Range("E34").Value = "  -0.67%  "

$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("D36").Value = "3.663.53"
$ws.Range("E36").Value = "  -0.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("E38").Value = "  +0.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.82"
$ws.Range("E39").Value = "  +0.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.132"
$ws.Range("E40").Value = "  -1.50%  "

$ws.Range("E41").Value = "  -0.84%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "436.45"
$ws.Range("E42").Value = "  -2.59%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "49.06"
$ws.Range("E43").Value = "  -0.29%  "

$ws.Range("E44").Value = "  -1.58%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.87"
$ws.Range("E45").Value = "  -1.85%  "

$ws.Range("E46").Value = "  +1.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.76"
$ws.Range("E47").Value = "  -1.22%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.03"
$ws.Range("E49").Value = "  +1.73%  "

$ws.Range("E50").Value = "  +0.68%  "

$ws.Range("D51").Value = "2.747.84"
$ws.Range("E51").Value = "  -2.88%  "
